# Refresh the cryptos worksheet with the latest scraped price/volume
# snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price strings (column D) look numeric ("65.440.32", "0.999", "6.78", ...)
# so a plain .Value assignment would make Excel silently coerce them to
# numbers (and mangle the dotted "thousands" values / lose trailing zeros).
# Forcing a text number-format while writing keeps them as literal strings,
# then restoring the "Normal" style afterwards so no stray cell format is
# left behind (matches the original / target workbook, which has no style
# attribute on these cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "65.440.32"
$ws.Range("E2").Value = "  +0.81%  "
Set-TextValue $ws.Range("D3") "2.958.54"
$ws.Range("E3").Value = "  -0.54%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "567.72"
$ws.Range("E5").Value = "  -2.28%  "
Set-TextValue $ws.Range("D6") "158.61"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue $ws.Range("D8") "0.521"
$ws.Range("E8").Value = "  +1.48%  "
Set-TextValue $ws.Range("D9") "2.955.38"
$ws.Range("E9").Value = "  -0.52%  "
Set-TextValue $ws.Range("D10") "6.78"
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("E11").Value = "  +0.76%  "
Set-TextValue $ws.Range("D12") "0.457"
$ws.Range("E12").Value = "  +2.60%  "
Set-TextValue $ws.Range("D13") "0.0000245"
$ws.Range("E13").Value = "  +3.50%  "
Set-TextValue $ws.Range("D14") "34.16"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("E15").Value = "  -0.51%  "
Set-TextValue $ws.Range("D16") "65.513.66"
$ws.Range("E16").Value = "  +1.01%  "
Set-TextValue $ws.Range("D17") "3.445.17"
$ws.Range("E17").Value = "  -0.66%  "
Set-TextValue $ws.Range("D18") "6.97"
$ws.Range("E18").Value = "  +1.65%  "
Set-TextValue $ws.Range("D19") "2.956.29"
$ws.Range("E19").Value = "  -0.63%  "
Set-TextValue $ws.Range("D20") "447.25"
$ws.Range("E20").Value = "  +0.06%  "
Set-TextValue $ws.Range("D21") "13.89"
$ws.Range("E21").Value = "  +2.14%  "
Set-TextValue $ws.Range("D22") "0.680"
$ws.Range("E22").Value = "  +0.54%  "
Set-TextValue $ws.Range("D23") "7.22"
$ws.Range("E23").Value = "  +0.03%  "
Set-TextValue $ws.Range("D24") "82.88"
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D25") "2.18"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D26") "12.08"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").Value = "  +0.03%  "
Set-TextValue $ws.Range("D28") "9.97"
$ws.Range("E28").Value = "  -5.40%  "
Set-TextValue $ws.Range("D29") "7.94"
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E30").Value = "  -0.87%  "
Set-TextValue $ws.Range("D31") "2.57"
$ws.Range("E31").Value = "  -0.18%  "
Set-TextValue $ws.Range("D32") "0.0₃0975"
$ws.Range("E32").Value = "  -3.64%  "
Set-TextValue $ws.Range("D33") "27.48"
$ws.Range("E33").Value = "  +3.71%  "
$ws.Range("E34").Value = "  +0.64%  "
Set-TextValue $ws.Range("D35") "0.998"
$ws.Range("E35").Value = "  -0.06%  "
Set-TextValue $ws.Range("D36") "0.975"
$ws.Range("E36").Value = "  -0.27%  "
Set-TextValue $ws.Range("D37") "5.74"
$ws.Range("E37").Value = "  +1.83%  "
Set-TextValue $ws.Range("D38") "49.11"
$ws.Range("E38").Value = "  +0.58%  "
Set-TextValue $ws.Range("D39") "1.98"
$ws.Range("E39").Value = "  -4.49%  "
Set-TextValue $ws.Range("D40") "0.299"
$ws.Range("E40").Value = "  +1.54%  "
Set-TextValue $ws.Range("D41") "0.119"
$ws.Range("E41").Value = "  -1.01%  "
Set-TextValue $ws.Range("D42") "42.94"
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D43") "2.78"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D44") "8.46"
$ws.Range("E44").Value = "  +1.00%  "
Set-TextValue $ws.Range("D45") "385.68"
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("E46").Value = "  +1.65%  "
Set-TextValue $ws.Range("D47") "2.740.66"
$ws.Range("E47").Value = "  -0.69%  "
Set-TextValue $ws.Range("D48") "130.71"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("E49").Value = "  +0.04%  "
Set-TextValue $ws.Range("D50") "0.107"
$ws.Range("E50").Value = "  +1.60%  "
Set-TextValue $ws.Range("D51") "2.15"
$ws.Range("E51").Value = "  +6.33%  "
